$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 0

# Row 11
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0

# Row 14
$ws.Range("L14").Value = 14
$ws.Range("M14").Value = 0

# Row 46
$ws.Range("L46").Value = 8
$ws.Range("M46").Value = 0

# Row 50
$ws.Range("L50").Value = 12
$ws.Range("M50").Value = 0

# Row 53
$ws.Range("L53").Value = 4
$ws.Range("M53").Value = 0

# Row 56
$ws.Range("L56").Value = 1
$ws.Range("M56").Value = 0

# Row 69
$ws.Range("L69").Value = 4
$ws.Range("M69").Value = 0

# Resumen / metricas de resumen section
# Total_Unidades
$ws.Range("C82").Value = 237

# Total_Ajuste_Stock
$ws.Range("C93").Value = 0
